$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"2.430374"
$ws.Cells.Item(2, 8).Value = [double]"7.291122000000001"
$ws.Cells.Item(2, 9).Value = [double]"0.009222757332915244"
$ws.Cells.Item(2, 10).Value = [double]"0.009222757332915246"
$ws.Cells.Item(2, 13).Value = [double]"168.1098273333333"
$ws.Cells.Item(2, 14).Value = [double]"504.329482"
$ws.Cells.Item(2, 15).Value = [double]"0.2984182258032519"
$ws.Cells.Item(2, 16).Value = [double]"0.298418225803252"
$ws.Cells.Item(2, 17).Value = [double]"408.5697534954227"
$ws.Cells.Item(2, 18).Value = [double]"3677.127781458804"
$ws.Cells.Item(2, 19).Value = [double]"0.002752238880302499"
$ws.Cells.Item(2, 20).Value = [double]"0.0027522388803025"
$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"2.430374"
$ws.Cells.Item(3, 8).Value = [double]"7.291122000000001"
$ws.Cells.Item(3, 9).Value = [double]"0.009222757332915244"
$ws.Cells.Item(3, 10).Value = [double]"0.009222757332915246"
$ws.Cells.Item(3, 15).Value = [double]"0.2893586437755394"
$ws.Cells.Item(3, 16).Value = [double]"0.2893586437755394"
$ws.Cells.Item(3, 17).Value = [double]"396.1661170021393"
$ws.Cells.Item(3, 18).Value = [double]"3565.495053019254"
$ws.Cells.Item(3, 19).Value = [double]"0.002668684553723266"
$ws.Cells.Item(3, 20).Value = [double]"0.002668684553723267"
$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"2.430374"
$ws.Cells.Item(4, 8).Value = [double]"7.291122000000001"
$ws.Cells.Item(4, 9).Value = [double]"0.009222757332915244"
$ws.Cells.Item(4, 10).Value = [double]"0.009222757332915246"
$ws.Cells.Item(4, 13).Value = [double]"165.99353"
$ws.Cells.Item(4, 14).Value = [double]"497.98059"
$ws.Cells.Item(4, 15).Value = [double]"0.294661504941043"
$ws.Cells.Item(4, 16).Value = [double]"0.294661504941043"
$ws.Cells.Item(4, 17).Value = [double]"403.42635948022"
$ws.Cells.Item(4, 18).Value = [double]"3630.83723532198"
$ws.Cells.Item(4, 19).Value = [double]"0.002717591555422845"
$ws.Cells.Item(4, 20).Value = [double]"0.002717591555422846"
$ws.Cells.Item(5, 5).Value = [double]"3"
$ws.Cells.Item(5, 6).Value = [double]"1"
$ws.Cells.Item(5, 7).Value = [double]"2.430374"
$ws.Cells.Item(5, 8).Value = [double]"7.291122000000001"
$ws.Cells.Item(5, 9).Value = [double]"0.009222757332915244"
$ws.Cells.Item(5, 10).Value = [double]"0.009222757332915246"
$ws.Cells.Item(5, 13).Value = [double]"66.22673433333334"
$ws.Cells.Item(5, 14).Value = [double]"198.680203"
$ws.Cells.Item(5, 15).Value = [double]"0.1175616254801657"
$ws.Cells.Item(5, 16).Value = [double]"0.1175616254801657"
$ws.Cells.Item(5, 17).Value = [double]"160.9557332286407"
$ws.Cells.Item(5, 18).Value = [double]"1448.601599057766"
$ws.Cells.Item(5, 19).Value = [double]"0.001084242343466633"
$ws.Cells.Item(5, 20).Value = [double]"0.001084242343466634"
$ws.Cells.Item(6, 9).Value = [double]"0.5480399755605952"
$ws.Cells.Item(6, 10).Value = [double]"0.5480399755605954"
$ws.Cells.Item(6, 13).Value = [double]"168.1098273333333"
$ws.Cells.Item(6, 14).Value = [double]"504.329482"
$ws.Cells.Item(6, 15).Value = [double]"0.2984182258032519"
$ws.Cells.Item(6, 16).Value = [double]"0.298418225803252"
$ws.Cells.Item(6, 17).Value = [double]"24278.26621018259"
$ws.Cells.Item(6, 18).Value = [double]"218504.3958916433"
$ws.Cells.Item(6, 19).Value = [double]"0.1635451171760504"
$ws.Cells.Item(6, 20).Value = [double]"0.1635451171760504"
$ws.Cells.Item(7, 9).Value = [double]"0.5480399755605952"
$ws.Cells.Item(7, 10).Value = [double]"0.5480399755605954"
$ws.Cells.Item(7, 15).Value = [double]"0.2893586437755394"
$ws.Cells.Item(7, 16).Value = [double]"0.2893586437755394"
$ws.Cells.Item(7, 17).Value = [double]"23541.2102088953"
$ws.Cells.Item(7, 18).Value = [double]"211870.8918800577"
$ws.Cells.Item(7, 19).Value = [double]"0.1585801040629936"
$ws.Cells.Item(7, 20).Value = [double]"0.1585801040629937"
$ws.Cells.Item(8, 9).Value = [double]"0.5480399755605952"
$ws.Cells.Item(8, 10).Value = [double]"0.5480399755605954"
$ws.Cells.Item(8, 13).Value = [double]"165.99353"
$ws.Cells.Item(8, 14).Value = [double]"497.98059"
$ws.Cells.Item(8, 15).Value = [double]"0.294661504941043"
$ws.Cells.Item(8, 16).Value = [double]"0.294661504941043"
$ws.Cells.Item(8, 17).Value = [double]"23972.63250123416"
$ws.Cells.Item(8, 18).Value = [double]"215753.6925111075"
$ws.Cells.Item(8, 19).Value = [double]"0.1614862839665374"
$ws.Cells.Item(8, 20).Value = [double]"0.1614862839665374"
$ws.Cells.Item(9, 9).Value = [double]"0.5480399755605952"
$ws.Cells.Item(9, 10).Value = [double]"0.5480399755605954"
$ws.Cells.Item(9, 13).Value = [double]"66.22673433333334"
$ws.Cells.Item(9, 14).Value = [double]"198.680203"
$ws.Cells.Item(9, 15).Value = [double]"0.1175616254801657"
$ws.Cells.Item(9, 16).Value = [double]"0.1175616254801657"
$ws.Cells.Item(9, 17).Value = [double]"9564.403889295367"
$ws.Cells.Item(9, 18).Value = [double]"86079.6350036583"
$ws.Cells.Item(9, 19).Value = [double]"0.06442847035501384"
$ws.Cells.Item(9, 20).Value = [double]"0.06442847035501385"
$ws.Cells.Item(10, 7).Value = [double]"116.470388"
$ws.Cells.Item(10, 8).Value = [double]"349.411164"
$ws.Cells.Item(10, 9).Value = [double]"0.4419805861132828"
$ws.Cells.Item(10, 10).Value = [double]"0.4419805861132828"
$ws.Cells.Item(10, 13).Value = [double]"168.1098273333333"
$ws.Cells.Item(10, 14).Value = [double]"504.329482"
$ws.Cells.Item(10, 15).Value = [double]"0.2984182258032519"
$ws.Cells.Item(10, 16).Value = [double]"0.298418225803252"
$ws.Cells.Item(10, 17).Value = [double]"19579.81681612634"
$ws.Cells.Item(10, 18).Value = [double]"176218.351345137"
$ws.Cells.Item(10, 19).Value = [double]"0.1318950623474073"
$ws.Cells.Item(10, 20).Value = [double]"0.1318950623474073"
$ws.Cells.Item(11, 7).Value = [double]"116.470388"
$ws.Cells.Item(11, 8).Value = [double]"349.411164"
$ws.Cells.Item(11, 9).Value = [double]"0.4419805861132828"
$ws.Cells.Item(11, 10).Value = [double]"0.4419805861132828"
$ws.Cells.Item(11, 15).Value = [double]"0.2893586437755394"
$ws.Cells.Item(11, 16).Value = [double]"0.2893586437755394"
$ws.Cells.Item(11, 17).Value = [double]"18985.39951451611"
$ws.Cells.Item(11, 18).Value = [double]"170868.5956306449"
$ws.Cells.Item(11, 19).Value = [double]"0.1278909029728575"
$ws.Cells.Item(11, 20).Value = [double]"0.1278909029728575"
$ws.Cells.Item(12, 7).Value = [double]"116.470388"
$ws.Cells.Item(12, 8).Value = [double]"349.411164"
$ws.Cells.Item(12, 9).Value = [double]"0.4419805861132828"
$ws.Cells.Item(12, 10).Value = [double]"0.4419805861132828"
$ws.Cells.Item(12, 13).Value = [double]"165.99353"
$ws.Cells.Item(12, 14).Value = [double]"497.98059"
$ws.Cells.Item(12, 15).Value = [double]"0.294661504941043"
$ws.Cells.Item(12, 16).Value = [double]"0.294661504941043"
$ws.Cells.Item(12, 17).Value = [double]"19333.33084458964"
$ws.Cells.Item(12, 18).Value = [double]"173999.9776013067"
$ws.Cells.Item(12, 19).Value = [double]"0.1302346646588642"
$ws.Cells.Item(12, 20).Value = [double]"0.1302346646588642"
$ws.Cells.Item(13, 7).Value = [double]"116.470388"
$ws.Cells.Item(13, 8).Value = [double]"349.411164"
$ws.Cells.Item(13, 9).Value = [double]"0.4419805861132828"
$ws.Cells.Item(13, 10).Value = [double]"0.4419805861132828"
$ws.Cells.Item(13, 13).Value = [double]"66.22673433333334"
$ws.Cells.Item(13, 14).Value = [double]"198.680203"
$ws.Cells.Item(13, 15).Value = [double]"0.1175616254801657"
$ws.Cells.Item(13, 16).Value = [double]"0.1175616254801657"
$ws.Cells.Item(13, 17).Value = [double]"7713.453443776256"
$ws.Cells.Item(13, 18).Value = [double]"69421.0809939863"
$ws.Cells.Item(13, 19).Value = [double]"0.05195995613415386"
$ws.Cells.Item(13, 20).Value = [double]"0.05195995613415386"
$ws.Cells.Item(14, 5).Value = [double]"3"
$ws.Cells.Item(14, 6).Value = [double]"1"
$ws.Cells.Item(14, 7).Value = [double]"0.1994"
$ws.Cells.Item(14, 8).Value = [double]"0.5982000000000001"
$ws.Cells.Item(14, 9).Value = [double]"0.0007566809932065188"
$ws.Cells.Item(14, 10).Value = [double]"0.0007566809932065189"
$ws.Cells.Item(14, 13).Value = [double]"168.1098273333333"
$ws.Cells.Item(14, 14).Value = [double]"504.329482"
$ws.Cells.Item(14, 15).Value = [double]"0.2984182258032519"
$ws.Cells.Item(14, 16).Value = [double]"0.298418225803252"
$ws.Cells.Item(14, 17).Value = [double]"33.52109957026667"
$ws.Cells.Item(14, 18).Value = [double]"301.6898961324"
$ws.Cells.Item(14, 19).Value = [double]"0.0002258073994917319"
$ws.Cells.Item(14, 20).Value = [double]"0.0002258073994917319"
$ws.Cells.Item(15, 5).Value = [double]"3"
$ws.Cells.Item(15, 6).Value = [double]"1"
$ws.Cells.Item(15, 7).Value = [double]"0.1994"
$ws.Cells.Item(15, 8).Value = [double]"0.5982000000000001"
$ws.Cells.Item(15, 9).Value = [double]"0.0007566809932065188"
$ws.Cells.Item(15, 10).Value = [double]"0.0007566809932065189"
$ws.Cells.Item(15, 15).Value = [double]"0.2893586437755394"
$ws.Cells.Item(15, 16).Value = [double]"0.2893586437755394"
$ws.Cells.Item(15, 17).Value = [double]"32.50344339193334"
$ws.Cells.Item(15, 18).Value = [double]"292.5309905274"
$ws.Cells.Item(15, 19).Value = [double]"0.0002189521859649664"
$ws.Cells.Item(15, 20).Value = [double]"0.0002189521859649665"
$ws.Cells.Item(16, 5).Value = [double]"3"
$ws.Cells.Item(16, 6).Value = [double]"1"
$ws.Cells.Item(16, 7).Value = [double]"0.1994"
$ws.Cells.Item(16, 8).Value = [double]"0.5982000000000001"
$ws.Cells.Item(16, 9).Value = [double]"0.0007566809932065188"
$ws.Cells.Item(16, 10).Value = [double]"0.0007566809932065189"
$ws.Cells.Item(16, 13).Value = [double]"165.99353"
$ws.Cells.Item(16, 14).Value = [double]"497.98059"
$ws.Cells.Item(16, 15).Value = [double]"0.294661504941043"
$ws.Cells.Item(16, 16).Value = [double]"0.294661504941043"
$ws.Cells.Item(16, 17).Value = [double]"33.099109882"
$ws.Cells.Item(16, 18).Value = [double]"297.8919889380001"
$ws.Cells.Item(16, 19).Value = [double]"0.000222964760218516"
$ws.Cells.Item(16, 20).Value = [double]"0.000222964760218516"
$ws.Cells.Item(17, 5).Value = [double]"3"
$ws.Cells.Item(17, 6).Value = [double]"1"
$ws.Cells.Item(17, 7).Value = [double]"0.1994"
$ws.Cells.Item(17, 8).Value = [double]"0.5982000000000001"
$ws.Cells.Item(17, 9).Value = [double]"0.0007566809932065188"
$ws.Cells.Item(17, 10).Value = [double]"0.0007566809932065189"
$ws.Cells.Item(17, 13).Value = [double]"66.22673433333334"
$ws.Cells.Item(17, 14).Value = [double]"198.680203"
$ws.Cells.Item(17, 15).Value = [double]"0.1175616254801657"
$ws.Cells.Item(17, 16).Value = [double]"0.1175616254801657"
$ws.Cells.Item(17, 17).Value = [double]"13.20561082606667"
$ws.Cells.Item(17, 18).Value = [double]"118.8504974346"
$ws.Cells.Item(17, 19).Value = [double]"8.895664753130454E-05"
$ws.Cells.Item(17, 20).Value = [double]"8.895664753130456E-05"
